$d = $word.ActiveDocument
$d.Content.Find.Execute("0.0%", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
